{"js": "// Locate the paragraph that starts the \"2.0)\" exercise block - this is the\n// first paragraph that gets renumbered to \"2.1)\" and is the anchor for the\n// whole block of new content appended by this commit.\nconst body = context.document.body;\nconst searchResults = body.search(\"2.0)\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst anchorRange = searchResults.items[0];\nconst startParagraphs = anchorRange.paragraphs;\nstartParagraphs.load(\"items\");\nawait context.sync();\nconst startParagraph = startParagraphs.items[0];\n\n// The replacement covers everything from the start of that paragraph through\n// to the end of the document body (the last three paragraphs of the old\n// text: \"2.0) a) ...\", \"b) ...\", and the dangling \"c) \" paragraph that held\n// the _GoBack bookmark).\nconst bodyParagraphs = body.paragraphs;\nbodyParagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = bodyParagraphs.items[bodyParagraphs.items.length - 1];\n\nconst target = startParagraph.getRange(\"Start\").expandTo(lastParagraph.getRange(\"End\"));\n\n// New paragraph content, expressed as WordprocessingML, inserted via\n// Range.insertOoxml (the Office.js analogue of the COM Range.InsertXML):\n//   * 3 blank paragraphs before the renumbered \"2.1)\" heading\n//   * \"2.0)\" -> \"2.1)\" with a lastRenderedPageBreak + the run split into\n//     \"2.1\" and \")\"\n//   * the existing \"b)\" paragraph unchanged\n//   * the existing \"c) \" paragraph with new text appended as its own run\n//   * new \"d)\", \"e)\", \"f)\" paragraphs\n//   * a whole new \"2.2)\" exercise with a)-d), the last of which keeps the\n//     _GoBack bookmark and gets an underlined paragraph mark\nconst bodyFragment =\n  \"<w:p/><w:p/><w:p/>\" +\n  '<w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:proofErr w:type=\"gramStart\"/><w:r><w:lastRenderedPageBreak/><w:t>2.1</w:t></w:r><w:r><w:t>)</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> a) John \u00e9 elegante e Kathy gosta de John.</w:t></w:r></w:p>' +\n  '<w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:r><w:t xml:space=\"preserve\">        b) Todo homem \u00e9 elegante.</w:t></w:r></w:p>' +\n  '<w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:r><w:t xml:space=\"preserve\">        c) </w:t></w:r><w:r><w:t>Toda mulher gosta de todo homem elegante.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">        d) Nem todo homem elegante gosta da Kathy.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">        e) Nem toda mulher que \u00e9 bonita gosta de todos os homens elegantes</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">        f) John gosta de toda mulher bonita.</w:t></w:r></w:p>' +\n  '<w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>2.2)</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> a) verdade, pois afirmando que tem um x </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>am</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> A(x) \u00e9 a mesma coisa que afirmar que para todo x que \u00acA(x).</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">        b) </w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">        c) verdade, pois afirmando que a um x tal que \u00ac</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>A(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>x) \u00e9 a mesma coisa que afirmar que para todo x, A(x).</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">        d)</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>';\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  bodyFragment +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that starts the \"2.0)\" exercise block - this is the\n# first paragraph that gets renumbered to \"2.1)\" and is the anchor for the\n# whole block of new content appended by this commit.\n$anchorRange = $d.Content\n$anchorRange.Find.Execute(\"2.0)\") | Out-Null\n$startPara = $anchorRange.Paragraphs(1)\n\n# The replacement covers everything from the start of that paragraph through\n# to the end of the document body (the last three paragraphs of the old\n# text: \"2.0) a) ...\", \"b) ...\", and the dangling \"c) \" paragraph that held\n# the _GoBack bookmark).\n$target = $d.Range($startPara.Range.Start, $d.Content.End)\n\n# New paragraph content, expressed as WordprocessingML, inserted via\n# Range.InsertXML (the COM analogue of Office.js Range.insertOoxml):\n#   * 3 blank paragraphs before the renumbered \"2.1)\" heading\n#   * \"2.0)\" -> \"2.1)\" with a lastRenderedPageBreak + the run split into\n#     \"2.1\" and \")\"\n#   * the existing \"b)\" paragraph unchanged\n#   * the existing \"c) \" paragraph with new text appended as its own run\n#   * new \"d)\", \"e)\", \"f)\" paragraphs\n#   * a whole new \"2.2)\" exercise with a)-d), the last of which keeps the\n#     _GoBack bookmark and gets an underlined paragraph mark\n$bodyFragment = @'\n<w:p/><w:p/><w:p/><w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:proofErr w:type=\"gramStart\"/><w:r><w:lastRenderedPageBreak/><w:t>2.1</w:t></w:r><w:r><w:t>)</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> a) John \u00e9 elegante e Kathy gosta de John.</w:t></w:r></w:p><w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:r><w:t xml:space=\"preserve\">        b) Todo homem \u00e9 elegante.</w:t></w:r></w:p><w:p w:rsidR=\"002C6299\" w:rsidRDefault=\"002C6299\" w:rsidP=\"009E0D23\"><w:r><w:t xml:space=\"preserve\">        c) </w:t></w:r><w:r><w:t>Toda mulher gosta de todo homem elegante.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">        d) Nem todo homem elegante gosta da Kathy.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">        e) Nem toda mulher que \u00e9 bonita gosta de todos os homens elegantes</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">        f) John gosta de toda mulher bonita.</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>2.2)</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> a) verdade, pois afirmando que tem um x </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>am</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> A(x) \u00e9 a mesma coisa que afirmar que para todo x que \u00acA(x).</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">        b) </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">        c) verdade, pois afirmando que a um x tal que \u00ac</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>A(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>x) \u00e9 a mesma coisa que afirmar que para todo x, A(x).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">        d)</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>\n'@\n\n$packageXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($packageXml)\n\nWrite-Output \"paragraphs now: $($d.Paragraphs.Count)\"\n"}
